$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 315, shifting rows 315:366 down to 316:367
$ws.Rows.Item(315).Insert()

# Populate the new row 315 with the latest week's data
$ws.Cells.Item(315, 1).Value = 8
$ws.Cells.Item(315, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(315, 3).Value = "Coquimbo"
$ws.Cells.Item(315, 4).Value = 45180
$ws.Cells.Item(315, 5).Value = 4
$ws.Cells.Item(315, 6).Value = 100112037
$ws.Cells.Item(315, 7).Value = "Cebollín"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 1100
$ws.Cells.Item(315, 11).Value = 1000
$ws.Cells.Item(315, 12).Value = 1200
$ws.Cells.Item(315, 13).Value = 1100
$ws.Cells.Item(315, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(315, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(315, 16).Value = 183
$ws.Cells.Item(315, 17).Value = 6
$ws.Cells.Item(315, 18).Value = "Hortaliza"
